$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.402.58'
$ws.Range("E2").Value = '  -2.32%  '
$ws.Range("D3").Value = '2.215.45'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '107.75'
$ws.Range("E5").Value = '  -12.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '294.72'
$ws.Range("E6").Value = '  +10.62%  '
$ws.Range("E7").Value = '  -3.74%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -4.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.61'
$ws.Range("E10").Value = '  -9.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").Value = '  -4.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.58'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.78'
$ws.Range("E13").Value = '  -4.86%  '
$ws.Range("E14").Value = '  -3.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.959'
$ws.Range("E15").Value = '  +5.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.93'
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("D17").Value = '2.548.06'
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").Value = '2.223.48'
$ws.Range("E18").Value = '  -1.98%  '
$ws.Range("D19").Value = '42.333.52'
$ws.Range("E19").Value = '  -2.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  +5.92%  '
$ws.Range("E21").Value = '  -4.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.56'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("E23").Value = '  +18.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("E24").Value = '  -6.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '227.68'
$ws.Range("E25").Value = '  -3.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.03'
$ws.Range("E26").Value = '  -4.73%  '
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.59'
$ws.Range("E28").Value = '  -3.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.96'
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.52'
$ws.Range("E30").Value = '  -9.62%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  -1.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("E32").Value = '  -4.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.60'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.86'
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0886'
$ws.Range("E35").Value = '  -3.20%  '
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.00'
$ws.Range("E37").Value = '  +8.25%  '
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("E39").Value = '  -3.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0363'
$ws.Range("E40").Value = '  -3.37%  '
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  -5.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.11'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("E44").Value = '  -3.92%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.69'
$ws.Range("E46").Value = '  -9.53%  '
$ws.Range("E47").Value = '  -5.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.42'
$ws.Range("E48").Value = '  -4.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.32'
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.18'
$ws.Range("E50").Value = '  +1.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.42'
$ws.Range("E51").Value = '  -1.74%  '
